$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure ambiguous numeric-looking price strings remain text (match source formatting, e.g. trailing zeros)
foreach ($addr in @("D5", "D7", "D9", "D10", "D11", "D13", "D15", "D16", "D19", "D21", "D23", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D35", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the crypto price refresh
$ws.Range("D2").Value = "37.129.37"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.070.77"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "252.82"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("D7").Value = "59.44"
$ws.Range("E7").Value = "  +10.16%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").Value = "61.53"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  +7.71%  "
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").Value = "16.27"
$ws.Range("E13").Value = "  +6.52%  "
$ws.Range("D14").Value = "2.373.99"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "0.819"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +7.67%  "
$ws.Range("D17").Value = "2.073.60"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "37.078.08"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "15.56"
$ws.Range("E19").Value = "  +8.00%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  +11.45%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "74.80"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("D23").Value = "240.66"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +13.50%  "
$ws.Range("D27").Value = "169.66"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "9.40"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "20.30"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("E31").Value = "  +5.27%  "
$ws.Range("D32").Value = "4.80"
$ws.Range("E32").Value = "  +6.53%  "
$ws.Range("D33").Value = "0.0637"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("E34").Value = "  +8.96%  "
$ws.Range("D35").Value = "0.0905"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").Value = "0.118"
$ws.Range("E38").Value = "  +29.76%  "
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").Value = "17.88"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "1.16"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "98.97"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "4.42"
$ws.Range("E45").Value = "  +16.05%  "
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("E47").Value = "  +12.86%  "
$ws.Range("E48").Value = "  +8.31%  "
$ws.Range("D49").Value = "1.304.79"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -1.53%  "
